$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (i.e. right
#    before the existing "总计" sheet).
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Copy the header-row formatting (bold/centered/bordered style) from the
# "2021-Q4" sheet so the new sheet's header looks the same.
$q4.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)

# Header labels.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row index cell (A2) formatted like the source sheet's index column.
$q4.Range("A2").Copy()
$q1.Range("A2").PasteSpecial(-4122)
$q1.Range("A2").Value = 0

# Data cells B2:G2 must stay plain text (matching the source workbook,
# which stores these numeric-looking values as text) without picking up
# any extra number-format style. Force text via NumberFormat, assign
# the value, then reset the cell format back to a plain/default style
# (copied from an untouched blank cell) so no style index is left on
# the cell -- exactly like the original sheets in this workbook.
$q1Data = @{
    "B2" = "162216"
    "C2" = "泰达宏利中证500指数增强（LOF）"
    "D2" = "4.44"
    "E2" = "93.67"
    "F2" = "1.26"
    "G2" = "0.0559"
}
foreach ($addr in $q1Data.Keys) {
    $cell = $q1.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $q1Data[$addr]
}
$q4.Range("AA1:AF1").Copy()
$q1.Range("B2:G2").PasteSpecial(-4122)

# H2 is a genuine number, same as the source sheets.
$q1.Range("H2").Value = 8

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: add a new first data row for
#    "2022-Q1", pushing the existing "2021-Q4" row down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()

# Re-stamp the A-column index style/values for both data rows: row 2
# is the new "2022-Q1" entry (index 0), row 3 is the shifted-down
# "2021-Q4" entry (index 1).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("A3").PasteSpecial(-4122)
$total.Range("A2").Value = 0
$total.Range("A3").Value = 1

$total.Range("B2").NumberFormat = "@"
$total.Range("B2").Value = "2022-Q1"

# Reset B2/C2/D2 to a plain/default style (no explicit style index),
# matching how the rest of the row data is stored in this workbook.
$q4.Range("AA1:AC1").Copy()
$total.Range("B2:D2").PasteSpecial(-4122)

$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.06

# Restore the original active sheet/selection (untouched by this edit).
$q4.Activate() | Out-Null
$q4.Range("A1").Select() | Out-Null
